$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-04-23"

# Update the label text for the April row (row 5)
$ws.Range("A5").Value = "April (through 04-23)"

# Revise March 2022 value (row 4, column I = year 2022)
$ws.Range("I4").Value = 133

# Update April 2022-row (row 5) data across all years as more data became available
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 25
$ws.Range("D5").Value = 45
$ws.Range("E5").Value = 40
$ws.Range("F5").Value = 35
$ws.Range("G5").Value = 48
$ws.Range("H5").Value = 85
$ws.Range("I5").Value = 103

# Update running Total row (row 6) across all years
$ws.Range("B6").Value = 80
$ws.Range("C6").Value = 153
$ws.Range("D6").Value = 234
$ws.Range("E6").Value = 237
$ws.Range("F6").Value = 145
$ws.Range("G6").Value = 246
$ws.Range("H6").Value = 508
$ws.Range("I6").Value = 538
